$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-4 (columns B and C) ---
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 22

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 3

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 4

# --- Add new rows 5-9 ---
$labels = @("3", "4", "5", "6", "7")
$winRates = @(100, 0, 100, 100, 100)
$gameLengths = @(6, 23, 12, 5, 9)

for ($i = 0; $i -lt 5; $i++) {
    $row = 5 + $i

    # Column A: copy the formatted "index" style from the row above it,
    # then force the cell to hold its number as literal text (matches
    # the existing A2:A4 label cells, which are text too).
    $ws.Range("A$($row - 1)").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)
    $ws.Range("A$row").NumberFormat = "@"
    $ws.Range("A$row").Value = $labels[$i]

    $ws.Range("B$row").Value = $winRates[$i]
    $ws.Range("C$row").Value = $gameLengths[$i]
}
